$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A7').Value = '2025-12-29 22:50:01'
$ws.Range('B7').Value = 'gemini-3-flash-preview'

$c7 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document f...
'@
$ws.Range('C7').Value = $c7

$d7 = @'

Role: You are a strict, pedantic Legal Proofreader. You are reviewing a standalone legal document fragment.
Input: The attached text from a contract.
CRITICAL INSTRUCTIONS:
1. **Assume Isolation with Common Sense**: Do NOT assume missing definitions exist in other documents. However, IGNORE common commercial lending terms typically defined in a base Credit Agreement (e.g., "Borrower", "Administrative Agent", "Lender", "Business Day", "Dollars", "GAAP", "Material Adverse Effect"). Only flag specific, deal-specific, or unusual capitalized terms that are undefined.
2. **Logic Check:** Check all math and logic tables.
3. **Drafting Errors:** Find any placeholders like "[__]" or blank lines that were forgotten.
Output Format:
Return ONLY a valid JSON object with the following structure:
{
  "errors": [
    {
      "location": "Page 3, Section 2.1",
      "error": "Description of the error",
      "suggestion": "Suggested fix"
    }
  ]
}
If no errors are found, return {"errors": []}.
--- CONTRACT TEXT BEGINS ---
--- [START OF PAGE 1] ---
FIRST AMENDMENT TO AMENDED AND RESTATED
CREDIT AGREEMENT
This document is a generated test file containing intentional legal drafting errors for AI training
purposes.
ARTICLE I: DEFINITIONS
...
"Applicable Margin" means the corresponding percentages per annum as set forth below based on
the Consolidated Total Leverage Ratio:
Pricing Level
Consolidated Total Leverage Ratio
SOFR Margin
I
Greater than or equal to 3.25 to 1.00
2.75%
II
Greater than or equal to 2.60 to 1.00 but less than 3.25 to 1.00
2.50%
III
Greater than or equal to 1.75 to 1.00 but less than 2.50 to 1.00
2.25%
IV
Greater than or equal to 1.50 to 1.00 but less than 1.75 to 1.00
2.00%
V
Less than 1.50 to 1.00
1.75%
"Cash Collateral" shall have a meaning correlative to the foregoing and shall include the proceeds of
such cash collateral...
(Note: "Cash Collateralize" is defined, but "Cash Collateralization" is NOT defined
in this section.)
SECTION 2.5 Permanent Reduction of the Revolving Credit
Commitment
The Borrower shall have the right at any time to terminate the Revolving Credit Commitment...
Any reduction of the Revolving Credit Commitment to zero shall be accompanied by payment of all
outstanding Revolving Credit Loans and furnishing of Cash Collateralization satisfactory to the
Administrative Agent.
SECTION 5.13 Incremental Loans
At any time after the First Amendment Effective Date, the Borrower may by written notice to the
Administrative Agent elect to request the establishment of one or more increases in the Revolving
--- [START OF PAGE 2] ---
Credit Commitments...
Such notice shall be delivered to the Administrative Agent at its office in [__] (or such other location as
the Administrative Agent may designate).
The terms of such Incremental Revolving Credit Increase shall be subject to the documentation
requirements set forth in Section 5.19.
[End of Document]
--- CONTRACT TEXT ENDS ---
'@
$ws.Range('D7').Value = $d7

$e7 = @'
{
  "errors": [
    {
      "location": "Page 1, Article I, 'Applicable Margin' table",
      "error": "Mathematical/logical gap in pricing levels: Level II ends at '2.60 to 1.00' and Level III starts at '2.50 to 1.00'. A ratio between 2.50 and 2.60 (e.g., 2.55) is not covered.",
      "suggestion": "Change Level III to 'Greater than or equal to 1.75 to 1.00 but less than 2.60 to 1.00'."
    },
    {
      "location": "Page 1, Section 2.5",
      "error": "The term 'Cash Collateralization' is used but is not defined in the Article I definitions section (as noted in the parenthetical within Article I).",
      "suggestion": "Add a definition for 'Cash Collateralization' or replace with the term 'Cash Collateral'."
    },
    {
      "location": "Page 2, Section 5.13",
      "error": "Drafting placeholder '[__]' remains in the text regarding the location of the Administrative Agent's office.",
      "suggestion": "Replace the placeholder with the appropriate city, state, or specific address."
    },
    {
      "location": "Page 2, Section 5.13",
      "error": "The term 'Incremental Revolving Credit Increase' is capitalized and used as a specific instrument/action but is not defined in this document.",
      "suggestion": "Define 'Incremental Revolving Credit Increase' in Article I or within Section 5.13."
    },
    {
      "location": "Page 2, Section 5.13",
      "error": "Internal cross-reference to 'Section 5.19' appears to be broken as the document fragment ends at Section 5.13.",
      "suggestion": "Verify if the documentation requirements should refer to a section within this Amendment or a section in the base Credit Agreement."
    }
  ]
}
'@
$ws.Range('E7').Value = $e7

Write-Output "done"